# Generate Report for Archive
#
# 1) The "Status" value shown for this handback entry moves from
#    "Ready for handoff" to "In Translation" on every sheet that surfaces
#    it (Overview!E2/F2, zh-cn!C2, de-de!C2 all pointed at the same shared
#    string, so updating every occurrence collapses back onto one string).
# 2) The Status column is narrowed (it was sized the same as a datetime
#    column; now it is sized like the other short columns).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Status text: "Ready for handoff" -> "In Translation" -------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status columns -----------------------------------------
# Target stored width ~13.41 chars; ColumnWidth is quantized to the
# on-screen pixel grid by the host, so 12.5 is the closest achievable
# setting (lands on the same pixel width Excel would show for ~13.41).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
